# Daily attendance processing - 2025-12-05 19:24:32
# Reverse the order of the comma-separated "Recorded By" entries in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ",\s*"
        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $cell.Value2 = [string]::Join(", ", $reversed)
        }
    }
}
